$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.818.66"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "2.103.87"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.58"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.08"
$ws.Range("E7").Value = "  +2.95%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.389"
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0841"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.76"
$ws.Range("E12").Value = "  +6.23%  "
$ws.Range("D13").Value = "2.416.25"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.99"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.806"
$ws.Range("E15").Value = "  +3.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.52"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").Value = "2.099.17"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").Value = "38.805.07"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.72"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.10"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("D21").Value = "0.0₃0848"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.00"
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.36"
$ws.Range("E24").Value = "  -2.84%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.73"
$ws.Range("E26").Value = "  +3.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.69"
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("E29").Value = "  +4.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.35"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.55"
$ws.Range("E31").Value = "  +9.94%  "
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.58"
$ws.Range("E33").Value = "  +2.15%  "
$ws.Range("B34").Value = "THORChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.23"
$ws.Range("E34").Value = "  +13.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.74"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0617"
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.38"
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.02"
$ws.Range("E40").Value = "  -1.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0227"
$ws.Range("E41").Value = "  +3.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.76"
$ws.Range("E42").Value = "  +1.58%  "
$ws.Range("D43").Value = "1.525.08"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("E44").Value = "  +8.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.80"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.81"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0915"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("E48").Value = "  +6.00%  "
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("D51").Value = "2.303.31"
$ws.Range("E51").Value = "  +1.11%  "
